$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the duplicated "Area / Nearby Station" values that were
# mistakenly copied into column H (rows 4-12), keeping formatting intact.
$ws.Range("H4:H12").ClearContents()

# Move the active selection from H8 to M2.
$ws.Range("M2").Select()
